$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.176.21'
$ws.Range("E2").Value = '  -3.97%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.690.85'
$ws.Range("E3").Value = '  -4.60%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.27'
$ws.Range("E5").Value = '  +0.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.82'
$ws.Range("E6").Value = '  +8.58%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.682.98'
$ws.Range("E7").Value = '  -4.58%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.627'
$ws.Range("E8").Value = '  -6.79%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.06%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.715'
$ws.Range("E10").Value = '  -5.18%  '

$ws.Range("E11").Value = '  -8.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '56.52'
$ws.Range("E12").Value = '  +5.71%  '

$ws.Range("E13").Value = '  -9.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.38'
$ws.Range("E14").Value = '  -9.47%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.271.97'
$ws.Range("E15").Value = '  -4.90%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.689.66'
$ws.Range("E16").Value = '  -4.68%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.30'
$ws.Range("E17").Value = '  -9.89%  '

$ws.Range("E18").Value = '  -2.34%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.82'
$ws.Range("E19").Value = '  -7.61%  '

$ws.Range("E20").Value = '  -7.86%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '68.020.92'
$ws.Range("E21").Value = '  -4.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '408.54'
$ws.Range("E22").Value = '  -6.68%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.57'
$ws.Range("E23").Value = '  -3.19%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.94'
$ws.Range("E24").Value = '  -5.74%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.02'
$ws.Range("E25").Value = '  -7.74%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.81'
$ws.Range("E26").Value = '  -8.03%  '

$ws.Range("E27").Value = '  -3.38%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.89'
$ws.Range("E28").Value = '  -4.56%  '

$ws.Range("E29").Value = '  +1.64%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.42'
$ws.Range("E30").Value = '  -9.60%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.82'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.24'
$ws.Range("E32").Value = '  -10.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.46'
$ws.Range("E33").Value = '  -8.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '43.33'
$ws.Range("E35").Value = '  -10.48%  '

$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '64.39'
$ws.Range("E36").Value = '  -7.93%  '

$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '602.01'
$ws.Range("E37").Value = '  -5.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0886'
$ws.Range("E38").Value = '  -10.18%  '

$ws.Range("E39").Value = '  +0.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.399'
$ws.Range("E40").Value = '  -6.15%  '

$ws.Range("E41").Value = '  +0.08%  '

$ws.Range("E42").Value = '  -7.28%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.80'
$ws.Range("E43").Value = '  +2.80%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.00'
$ws.Range("E44").Value = '  -9.62%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0436'
$ws.Range("E45").Value = '  -7.38%  '

$ws.Range("E46").Value = '  -12.75%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.20'
$ws.Range("E47").Value = '  -8.91%  '

$ws.Range("E48").Value = '  -4.35%  '

$ws.Range("E49").Value = '  -6.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.750.94'
$ws.Range("E50").Value = '  -2.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.16'
$ws.Range("E51").Value = '  -4.10%  '
